$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F8").Formula = "=C8-C5"
$ws.Range("F15").Formula = "=C15-C12"
$ws.Range("F22").Formula = "=C22-C19"
$ws.Range("F29").Formula = "=C29-C26"
$ws.Range("F36").Formula = "=C36-C33"
$ws.Range("F43").Formula = "=C43-C40"
$ws.Range("F46").Formula = "=AVERAGE(F2:F43)"
$ws.Range("G45").Select() | Out-Null
